# Update the "diagnosedCasesList" table on Sheet1 per feedback from the
# ORTHO & PT unit member:
#   - remove the obsolete "Lumbar disorder / HNP" row
#   - add two new "Shoulder impingement.../Tendinitis" and
#     ".../impingement" rows after "Calcific Tendinitis"
#   - add a duplicate "Shoulder impingement.../adhesive capsulitis" row
#     just before the CVA / CP rows at the bottom of the list
#
# Net effect: the table grows from A1:B88 to A1:B90.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Remove the "Lumbar disorder | HNP" row (old row 60).
$ws.Rows.Item(60).Delete()

# 2) After the deletion, "Shoulder impingement/.../Calcific Tendinitis"
#    sits at row 84, immediately followed by "CVA" at row 85. Insert three
#    blank rows there: two for the new sub-classifications, one more that
#    duplicates the "adhesive capsulitis" sub-classification.
$ws.Rows.Item(85).Insert()
$ws.Rows.Item(86).Insert()
$ws.Rows.Item(87).Insert()

$shoulderClassification = "Shoulder impingement/rotator cuff tendinitis/adhesive capsulitis"

$ws.Cells.Item(85, 1).Value = $shoulderClassification
$ws.Cells.Item(85, 2).Value = "Tendinitis"

$ws.Cells.Item(86, 1).Value = $shoulderClassification
$ws.Cells.Item(86, 2).Value = "impingement"

$ws.Cells.Item(87, 1).Value = $shoulderClassification
$ws.Cells.Item(87, 2).Value = "adhesive capsulitis"

# 3) The named range "diagnosedCasesList" must grow to cover the two
#    extra rows (A1:B88 -> A1:B90).
$ws.Names.Item("diagnosedCasesList").RefersTo = "=Sheet1!`$A`$1:`$B`$90"

# 4) Restore the view state recorded at save time (scroll position /
#    active selection) to match the author's last working position.
$excel.ActiveWindow.ScrollRow = 70
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D82").Select()
